# Update the subtitle on slide 1 (version bump + revised date) prior to the
# JGF meeting.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$subtitle = $s.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange

# Paragraph 1: "Integration Overview, v1.1-draft" -> "Integration Overview, v1.0-draft"
# Reassign via a Characters() range spanning the whole paragraph so the two
# existing runs collapse into a single run, matching the target markup.
$para1 = $tr.Paragraphs(1)
$tr.Characters($para1.Start, $para1.Length).Text = "Integration Overview, v1.0-draft"

# Paragraph 2: "Last revised, 5/12/2016, ESK" -> "Last revised, 5/9/2016, ESK"
$para2 = $tr.Paragraphs(2)
$tr.Characters($para2.Start, $para2.Length).Text = "Last revised, 5/9/2016, ESK"
